$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.266759872436523
$ws.Range("B1").Value = 2.494807720184326
$ws.Range("C1").Value = 3.584692001342773
$ws.Range("D1").Value = 2.961567640304565
$ws.Range("E1").Value = 1.071389675140381
